# "Implementando classe de despesa"
# Adds the remaining Tag/* entries, a full Despesa/* CRUD block, and the
# corresponding batch of new error codes (E0026-E0040) to the API error
# catalogue sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B was missing for the last few pre-existing rows (Tag/Get,
# Tag/Pesquisar, Despesa/Create, Despesa/Update belong to rows whose
# Código/A column was already filled in).
$ws.Range("B23").Value = "Tag/Get"
$ws.Range("B24").Value = "Tag/Pesquisar"
$ws.Range("B25").Value = "Despesa/Create"
$ws.Range("B26").Value = "Despesa/Update"

# New error codes, rows 27-41 - reuse column A's existing text-number-format
# so the new cells keep the same look as the rest of the Código column.
$ws.Range("A26").Copy()
$ws.Range("A27:A41").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A27").Value = "E0026"
$ws.Range("A28").Value = "E0027"
$ws.Range("A29").Value = "E0028"
$ws.Range("A30").Value = "E0029"
$ws.Range("A31").Value = "E0030"
$ws.Range("A32").Value = "E0031"
$ws.Range("A33").Value = "E0032"
$ws.Range("A34").Value = "E0033"
$ws.Range("A35").Value = "E0034"
$ws.Range("A36").Value = "E0035"
$ws.Range("A37").Value = "E0036"
$ws.Range("A38").Value = "E0037"
$ws.Range("A39").Value = "E0038"
$ws.Range("A40").Value = "E0039"
$ws.Range("A41").Value = "E0040"

# Finish the Despesa CRUD block started above (Create/Update already set).
$ws.Range("B27").Value = "Despesa/Delete"
$ws.Range("B28").Value = "Despesa/Get"

# Leave the selection where the author left it.
$ws.Range("B28").Select()
